$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.932.64'
$ws.Range("E2").Value = '  +0.58%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.810.72'
$ws.Range("E3").Value = '  +1.60%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.26'
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4977'
$ws.Range("E7").Value = '  -2.58%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3908'
$ws.Range("E8").Value = '  +2.29%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09709'
$ws.Range("E9").Value = '  +24.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.099'
$ws.Range("E10").Value = '  +1.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '40.84'
$ws.Range("E11").Value = '  +0.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.420'
$ws.Range("E12").Value = '  +3.73%  '
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.46'
$ws.Range("E13").Value = '  +1.77%  '
$ws.Range("B14").Value = 'BinanceUSD'
$ws.Range("C14").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.001'
$ws.Range("E14").Value = '  -0.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.813.79'
$ws.Range("E15").Value = '  +1.81%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.279'
$ws.Range("E16").Value = '  +1.32%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001131'
$ws.Range("E17").Value = '  +5.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '92.18'
$ws.Range("E18").Value = '  +0.93%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06640'
$ws.Range("E19").Value = '  +1.42%  '
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.12'
$ws.Range("E21").Value = '  +0.76%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.909'
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.001.89'
$ws.Range("E23").Value = '  +0.65%  '
$ws.Range("E24").Value = '  +0.95%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.243'
$ws.Range("E25").Value = '  +0.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.89'
$ws.Range("E26").Value = '  -0.64%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.019.84'
$ws.Range("E27").Value = '  +1.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.52'
$ws.Range("E28").Value = '  +1.68%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.381'
$ws.Range("E29").Value = '  +1.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.55'
$ws.Range("E30").Value = '  +3.07%  '
$ws.Range("E31").Value = '  -0.80%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.033'
$ws.Range("E32").Value = '  +0.13%  '
$ws.Range("E33").Value = '  +1.35%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.608'
$ws.Range("E34").Value = '  -0.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.06721'
$ws.Range("E35").Value = '  -4.95%  '
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02325'
$ws.Range("E36").Value = '  +1.05%  '
$ws.Range("B37").Value = 'FraxShare'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.860'
$ws.Range("E37").Value = '  +0.50%  '
$ws.Range("E38").Value = '  +0.77%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.945'
$ws.Range("E39").Value = '  -0.81%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.23'
$ws.Range("E40").Value = '  -2.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6160'
$ws.Range("E41").Value = '  +1.26%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.163'
$ws.Range("E42").Value = '  +0.97%  '
$ws.Range("B43").Value = 'Frax'
$ws.Range("C43").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.13'
$ws.Range("E44").Value = '  +0.19%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5877'
$ws.Range("E45").Value = '  +0.09%  '
$ws.Range("B46").Value = 'WEMIXTOKEN'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.286'
$ws.Range("E46").Value = '  -2.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.687'
$ws.Range("E47").Value = '  -0.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '123.76'
$ws.Range("E48").Value = '  -1.62%  '
$ws.Range("E49").Value = '  +1.91%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.177'
$ws.Range("E50").Value = '  -1.48%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06758'
$ws.Range("E51").Value = '  -1.53%  '
